$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate so the old A:D layout doesn't leave stray
# formatting behind once everything shifts one column to the right.
$ws.Cells.Clear()

# ---- column widths ---------------------------------------------------
# New narrow "section" column A, plus the two content columns (C/D, was
# B/C before the shift) growing a bit wider.
$ws.Columns.Item(1).ColumnWidth = 13.1666666666667
$ws.Columns.Item(2).ColumnWidth = 26.7369791666667
$ws.Columns.Item(3).ColumnWidth = 47.5
$ws.Columns.Item(4).ColumnWidth = 39.3333333333333

# ---- header row --------------------------------------------------------
$ws.Range("B1").Value = "als"
$ws.Range("C1").Value = "kan ik"
$ws.Range("D1").Value = "zodat ik (resultaat)"
$ws.Range("E1").Value = "prioriteit volgens MoSCoW"
$ws.Range("B1:E1").Font.Bold = $true

# ---- section: Login Pagina ---------------------------------------------
$ws.Range("A2").Value = "Login Pagina"
$ws.Range("A2").Font.Bold = $true

$ws.Range("B3").Value = "Gebruiker"
$ws.Range("C3").Value = "Registreren"
$ws.Range("D3").Value = "Later kan inloggen en de website kan gebruiken"
$ws.Range("E3").Value = "Must"
$ws.Range("E3").Characters(1, 1).Font.Bold = $true

$ws.Range("B4").Value = "Gebruiker"
$ws.Range("C4").Value = "Inloggen"
$ws.Range("D4").Value = "De website kan gebruiken"
$ws.Range("E4").Value = "Must"
$ws.Range("E4").Font.Bold = $true

# ---- section: The Wall --------------------------------------------------
$ws.Range("A6").Value = "The Wall"
$ws.Range("A6").Font.Bold = $true

$ws.Range("B7").Value = "Gebruiker"
$ws.Range("C7").Value = "Eigen foto's uploaden"
$ws.Range("D7").Value = "The Wall kan vullen"
$ws.Range("E7").Value = "Must"
$ws.Range("E7").Font.Bold = $true

$ws.Range("B8").Value = "Gebruiker"
$ws.Range("C8").Value = "Contact pagina bezoeken"
$ws.Range("D8").Value = "De developers hun informatie kan vinden"
$ws.Range("E8").Value = "Could"
$ws.Range("E8").Font.Bold = $true

$ws.Range("B9").Value = "Gebruiker"
$ws.Range("C9").Value = "Zoeken naar trefwoorden"
$ws.Range("D9").Value = "Om een bepaalde foto te vinden"
$ws.Range("E9").Value = "Could"

$ws.Range("B10").Value = "Gebruiker "
$ws.Range("C10").Value = "Catogoriën als favoriet instellen"
$ws.Range("D10").Value = "Sneller naar mijn favoriete catogorien kan gaan"
$ws.Range("E10").Value = "Could"
$ws.Range("E10").Font.Bold = $true

$ws.Range("B11").Value = "Gebruiker"
$ws.Range("C11").Value = "Op categoriën kan klikken die ik zelf heb ingesteld"
$ws.Range("D11").Value = "Makkelijker bepaalde producten kan opzoeken"
$ws.Range("E11").Value = "Could"
$ws.Range("E11").Font.Bold = $true

# ---- section: Eigen Profiel ---------------------------------------------
$ws.Range("A12").Value = "Eigen Profiel"
$ws.Range("A12").Font.Bold = $true

$ws.Range("B13").Value = "Gebruiker"
$ws.Range("C13").Value = "Naar mijn eigen profiel"
$ws.Range("D13").Value = "Mijn geuploade foto's kan zien"
$ws.Range("E13").Value = "Could"
$ws.Range("E13").Font.Bold = $true

$ws.Range("B14").Value = "Gebruiker "
$ws.Range("C14").Value = "Naar mijn eigen profiel"
$ws.Range("D14").Value = "Kan uitloggen"
$ws.Range("E14").Value = "Could"
$ws.Range("E14").Font.Bold = $true

$ws.Range("B15").Value = "Gebruiker"
$ws.Range("C15").Value = "Naar mijn eigen profiel"
$ws.Range("D15").Value = "Om mijn profiel foto groter te zien"
$ws.Range("E15").Value = "Could"
$ws.Range("E15").Font.Bold = $true

# ---- section: Contact Pagina ---------------------------------------------
$ws.Range("A16").Value = "Contact Pagina"
$ws.Range("A16").Font.Bold = $true

$ws.Range("B17").Value = "Gebruiker "
$ws.Range("C17").Value = "Zien wat de email van de developers zijn"
$ws.Range("D17").Value = "Ze kan contacteren voor klachten/tips en tops"
$ws.Range("E17").Value = "Could"
$ws.Range("E17").Font.Bold = $true

$ws.Range("B18").Value = "Gebruiker"
$ws.Range("C18").Value = "Zien wat de telefoonnummers van de developers zijn"
$ws.Range("D18").Value = "Ze kan bellen voor vragen"
$ws.Range("E18").Value = "Could"
$ws.Range("E18").Font.Bold = $true

$ws.Range("B19").Value = "Gebruiker"
$ws.Range("C19").Value = "Zien wat de instagram links van de developers zijn"
$ws.Range("D19").Value = "Hun producten kan bekijken"
$ws.Range("E19").Value = "Could"
$ws.Range("E19").Font.Bold = $true

$ws.Range("B20").Value = "Gebruiker"
$ws.Range("C20").Value = "Op de link van de portfolio site van de developers klikken"
$ws.Range("D20").Value = "Hun projecten kan bekijken"
$ws.Range("E20").Value = "Could"
$ws.Range("E20").Font.Bold = $true

# ---- selection/view state ------------------------------------------------
$ws.Range("B27").Select()
